# Apply Wnt9a-Fzd4 NATMI update (Dr Hou advice): expand Target cluster
# categories from {ECs, FAPs, sCs} to {ECs, FAPs, M2, sCs} and refresh
# all computed metrics for the 3x4 sending/target cluster grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = New-Object 'object[,]' 12,20

# Row 2: ECs -> ECs
$newData[0,0] = "ECs"
$newData[0,1] = "Wnt9a"
$newData[0,2] = "Fzd4"
$newData[0,3] = "ECs"
$newData[0,4] = 3
$newData[0,5] = 1
$newData[0,6] = 0.7183383333333334
$newData[0,7] = 2.155015
$newData[0,8] = 0.1157638296793402
$newData[0,9] = 0.1157638296793402
$newData[0,10] = 3
$newData[0,11] = 1
$newData[0,12] = 30.38232766666667
$newData[0,13] = 91.146983
$newData[0,14] = 0.4410933377331532
$newData[0,15] = 0.4410933377331531
$newData[0,16] = 21.82479061886056
$newData[0,17] = 196.423115569745
$newData[0,18] = 0.05106265402203242
$newData[0,19] = 0.05106265402203242

# Row 3: ECs -> FAPs
$newData[1,0] = "ECs"
$newData[1,1] = "Wnt9a"
$newData[1,2] = "Fzd4"
$newData[1,3] = "FAPs"
$newData[1,4] = 3
$newData[1,5] = 1
$newData[1,6] = 0.7183383333333334
$newData[1,7] = 2.155015
$newData[1,8] = 0.1157638296793402
$newData[1,9] = 0.1157638296793402
$newData[1,10] = 3
$newData[1,11] = 1
$newData[1,12] = 22.690535
$newData[1,13] = 68.071605
$newData[1,14] = 0.3294232070665772
$newData[1,15] = 0.3294232070665772
$newData[1,16] = 16.29948109434167
$newData[1,17] = 146.695329849075
$newData[1,18] = 0.03813529203527726
$newData[1,19] = 0.03813529203527726

# Row 4: ECs -> M2
$newData[2,0] = "ECs"
$newData[2,1] = "Wnt9a"
$newData[2,2] = "Fzd4"
$newData[2,3] = "M2"
$newData[2,4] = 3
$newData[2,5] = 1
$newData[2,6] = 0.7183383333333334
$newData[2,7] = 2.155015
$newData[2,8] = 0.1157638296793402
$newData[2,9] = 0.1157638296793402
$newData[2,10] = 1
$newData[2,11] = 0.3333333333333333
$newData[2,12] = 0.09645133333333333
$newData[2,13] = 0.289354
$newData[2,14] = 0.001400289043537939
$newData[2,15] = 0.001400289043537939
$newData[2,16] = 0.06928469003444446
$newData[2,17] = 0.62356221031
$newData[2,18] = 0.0001621028223379722
$newData[2,19] = 0.0001621028223379721

# Row 5: ECs -> sCs
$newData[3,0] = "ECs"
$newData[3,1] = "Wnt9a"
$newData[3,2] = "Fzd4"
$newData[3,3] = "sCs"
$newData[3,4] = 3
$newData[3,5] = 1
$newData[3,6] = 0.7183383333333334
$newData[3,7] = 2.155015
$newData[3,8] = 0.1157638296793402
$newData[3,9] = 0.1157638296793402
$newData[3,10] = 3
$newData[3,11] = 1
$newData[3,12] = 15.71027466666667
$newData[3,13] = 47.130824
$newData[3,14] = 0.2280831661567317
$newData[3,15] = 0.2280831661567317
$newData[3,16] = 11.28529252026222
$newData[3,17] = 101.56763268236
$newData[3,18] = 0.02640378079969254
$newData[3,19] = 0.02640378079969253

# Row 6: FAPs -> ECs
$newData[4,0] = "FAPs"
$newData[4,1] = "Wnt9a"
$newData[4,2] = "Fzd4"
$newData[4,3] = "ECs"
$newData[4,4] = 3
$newData[4,5] = 1
$newData[4,6] = 4.267456
$newData[4,7] = 12.802368
$newData[4,8] = 0.687721964183189
$newData[4,9] = 0.687721964183189
$newData[4,10] = 3
$newData[4,11] = 1
$newData[4,12] = 30.38232766666667
$newData[4,13] = 91.146983
$newData[4,14] = 0.4410933377331532
$newData[4,15] = 0.4410933377331531
$newData[4,16] = 129.6552464950827
$newData[4,17] = 1166.897218455744
$newData[4,18] = 0.3033495766139628
$newData[4,19] = 0.3033495766139628

# Row 7: FAPs -> FAPs
$newData[5,0] = "FAPs"
$newData[5,1] = "Wnt9a"
$newData[5,2] = "Fzd4"
$newData[5,3] = "FAPs"
$newData[5,4] = 3
$newData[5,5] = 1
$newData[5,6] = 4.267456
$newData[5,7] = 12.802368
$newData[5,8] = 0.687721964183189
$newData[5,9] = 0.687721964183189
$newData[5,10] = 3
$newData[5,11] = 1
$newData[5,12] = 22.690535
$newData[5,13] = 68.071605
$newData[5,14] = 0.3294232070665772
$newData[5,15] = 0.3294232070665772
$newData[5,16] = 96.83085972896001
$newData[5,17] = 871.4777375606402
$newData[5,18] = 0.2265515750113519
$newData[5,19] = 0.2265515750113519

# Row 8: FAPs -> M2
$newData[6,0] = "FAPs"
$newData[6,1] = "Wnt9a"
$newData[6,2] = "Fzd4"
$newData[6,3] = "M2"
$newData[6,4] = 3
$newData[6,5] = 1
$newData[6,6] = 4.267456
$newData[6,7] = 12.802368
$newData[6,8] = 0.687721964183189
$newData[6,9] = 0.687721964183189
$newData[6,10] = 1
$newData[6,11] = 0.3333333333333333
$newData[6,12] = 0.09645133333333333
$newData[6,13] = 0.289354
$newData[6,14] = 0.001400289043537939
$newData[6,15] = 0.001400289043537939
$newData[6,16] = 0.4116018211413334
$newData[6,17] = 3.704416390272
$newData[6,18] = 0.0009630095314461106
$newData[6,19] = 0.0009630095314461104

# Row 9: FAPs -> sCs
$newData[7,0] = "FAPs"
$newData[7,1] = "Wnt9a"
$newData[7,2] = "Fzd4"
$newData[7,3] = "sCs"
$newData[7,4] = 3
$newData[7,5] = 1
$newData[7,6] = 4.267456
$newData[7,7] = 12.802368
$newData[7,8] = 0.687721964183189
$newData[7,9] = 0.687721964183189
$newData[7,10] = 3
$newData[7,11] = 1
$newData[7,12] = 15.71027466666667
$newData[7,13] = 47.130824
$newData[7,14] = 0.2280831661567317
$newData[7,15] = 0.2280831661567317
$newData[7,16] = 67.04290588791466
$newData[7,17] = 603.386152991232
$newData[7,18] = 0.1568578030264282
$newData[7,19] = 0.1568578030264282

# Row 10: sCs -> ECs
$newData[8,0] = "sCs"
$newData[8,1] = "Wnt9a"
$newData[8,2] = "Fzd4"
$newData[8,3] = "ECs"
$newData[8,4] = 3
$newData[8,5] = 1
$newData[8,6] = 1.219411
$newData[8,7] = 3.658233
$newData[8,8] = 0.1965142061374708
$newData[8,9] = 0.1965142061374708
$newData[8,10] = 3
$newData[8,11] = 1
$newData[8,12] = 30.38232766666667
$newData[8,13] = 91.146983
$newData[8,14] = 0.4410933377331532
$newData[8,15] = 0.4410933377331531
$newData[8,16] = 37.04854456233767
$newData[8,17] = 333.436901061039
$newData[8,18] = 0.0866811070971579
$newData[8,19] = 0.08668110709715787

# Row 11: sCs -> FAPs
$newData[9,0] = "sCs"
$newData[9,1] = "Wnt9a"
$newData[9,2] = "Fzd4"
$newData[9,3] = "FAPs"
$newData[9,4] = 3
$newData[9,5] = 1
$newData[9,6] = 1.219411
$newData[9,7] = 3.658233
$newData[9,8] = 0.1965142061374708
$newData[9,9] = 0.1965142061374708
$newData[9,10] = 3
$newData[9,11] = 1
$newData[9,12] = 22.690535
$newData[9,13] = 68.071605
$newData[9,14] = 0.3294232070665772
$newData[9,15] = 0.3294232070665772
$newData[9,16] = 27.669087974885
$newData[9,17] = 249.021791773965
$newData[9,18] = 0.06473634001994809
$newData[9,19] = 0.06473634001994807

# Row 12: sCs -> M2
$newData[10,0] = "sCs"
$newData[10,1] = "Wnt9a"
$newData[10,2] = "Fzd4"
$newData[10,3] = "M2"
$newData[10,4] = 3
$newData[10,5] = 1
$newData[10,6] = 1.219411
$newData[10,7] = 3.658233
$newData[10,8] = 0.1965142061374708
$newData[10,9] = 0.1965142061374708
$newData[10,10] = 1
$newData[10,11] = 0.3333333333333333
$newData[10,12] = 0.09645133333333333
$newData[10,13] = 0.289354
$newData[10,14] = 0.001400289043537939
$newData[10,15] = 0.001400289043537939
$newData[10,16] = 0.1176138168313333
$newData[10,17] = 1.058524351482
$newData[10,18] = 0.0002751766897538564
$newData[10,19] = 0.0002751766897538563

# Row 13: sCs -> sCs
$newData[11,0] = "sCs"
$newData[11,1] = "Wnt9a"
$newData[11,2] = "Fzd4"
$newData[11,3] = "sCs"
$newData[11,4] = 3
$newData[11,5] = 1
$newData[11,6] = 1.219411
$newData[11,7] = 3.658233
$newData[11,8] = 0.1965142061374708
$newData[11,9] = 0.1965142061374708
$newData[11,10] = 3
$newData[11,11] = 1
$newData[11,12] = 15.71027466666667
$newData[11,13] = 47.130824
$newData[11,14] = 0.2280831661567317
$newData[11,15] = 0.2280831661567317
$newData[11,16] = 19.15728174155467
$newData[11,17] = 172.415535673992
$newData[11,18] = 0.04482158233061097
$newData[11,19] = 0.04482158233061097

$ws.Range("A2:T13").Value = $newData
